# Update "想去人数" (F column) figures across the sheets, matching the
# newer scrape snapshot recorded in the commit.

$wb = $excel.ActiveWorkbook

# Map: sheet name -> @{ row = newValue }
$updates = @{
    "展览" = @{
        4  = 608
        7  = 5487
        8  = 1533
        10 = 3139
        13 = 1301
        14 = 4350
        16 = 1678
        17 = 2609
        18 = 39
        19 = 36
        22 = 993
        24 = 78
        26 = 86
        28 = 1098
        29 = 393
        30 = 61
        31 = 180
        32 = 312
        35 = 1693
        36 = 2209
        37 = 1028
        38 = 39
        39 = 256
        40 = 615
        41 = 312
        43 = 8
        44 = 659
        45 = 11
        46 = 409
        47 = 345
        48 = 213
        49 = 141
    }
    "演出" = @{
        10 = 151
    }
    "本地生活" = @{
        2 = 753
    }
    "全部类型" = @{
        2  = 753
        6  = 608
        9  = 5487
        10 = 1533
        13 = 3139
        15 = 1301
        16 = 4350
        18 = 1678
        20 = 39
        23 = 36
        25 = 151
        26 = 993
        28 = 86
        30 = 1098
        31 = 393
        32 = 61
        33 = 180
        35 = 1693
        36 = 2209
        37 = 1028
        38 = 39
        41 = 256
        42 = 615
        43 = 312
        44 = 659
        45 = 409
        46 = 345
        47 = 213
        48 = 141
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}
